$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-7 (existing rows) and new rows 8-10
$data = @{
    2  = @("ECs", "Efna1", "Epha3", "ECs", 3, 1, 23.60223533333333, 70.80670599999999, 0.8824726436021215, 0.8824726436021214, 2, 0.6666666666666666, 1.106217, 3.318651, 0.01813136626967656, 0.01813136626967656, 26.109193963734, 234.982745673606, 0.01600043472411981, 0.01600043472411981)
    3  = @("ECs", "Efna1", "Epha3", "FAPs", 3, 1, 23.60223533333333, 70.80670599999999, 0.8824726436021215, 0.8824726436021214, 3, 1, 59.24481466666666, 177.734444, 0.971047664518299, 0.971047664518299, 1398.310058042385, 12584.79052238146, 0.8569229995711294, 0.8569229995711293)
    4  = @("ECs", "Efna1", "Epha3", "sCs", 3, 1, 23.60223533333333, 70.80670599999999, 0.8824726436021215, 0.8824726436021214, 3, 1, 0.6602006666666667, 1.980602, 0.01082096921202439, 0.01082096921202439, 15.58221150189022, 140.239903517012, 0.009549209306872326, 0.009549209306872324)
    5  = @("FAPs", "Efna1", "Epha3", "ECs", 3, 1, 2.633202, 7.899606, 0.09845375648791208, 0.09845375648791205, 2, 0.6666666666666666, 1.106217, 3.318651, 0.01813136626967656, 0.01813136626967656, 2.912892816834, 26.216035351506, 0.001785101119507879, 0.001785101119507879)
    6  = @("FAPs", "Efna1", "Epha3", "FAPs", 3, 1, 2.633202, 7.899606, 0.09845375648791208, 0.09845375648791205, 3, 1, 59.24481466666666, 177.734444, 0.971047664518299, 0.971047664518299, 156.003564469896, 1404.032080229064, 0.09560329030064035, 0.09560329030064033)
    7  = @("FAPs", "Efna1", "Epha3", "sCs", 3, 1, 2.633202, 7.899606, 0.09845375648791208, 0.09845375648791205, 3, 1, 0.6602006666666667, 1.980602, 0.01082096921202439, 0.01082096921202439, 1.738441715868, 15.645975442812, 0.001065365067763843, 0.001065365067763842)
    8  = @("sCs", "Efna1", "Epha3", "ECs", 3, 1, 0.5101343333333334, 1.530403, 0.01907359990996641, 0.0190735999099664, 2, 0.6666666666666666, 1.106217, 3.318651, 0.01813136626967656, 0.01813136626967656, 0.564319271817, 5.078873446353, 0.0003458304260488708, 0.0003458304260488707)
    9  = @("sCs", "Efna1", "Epha3", "FAPs", 3, 1, 0.5101343333333334, 1.530403, 0.01907359990996641, 0.0190735999099664, 3, 1, 59.24481466666666, 177.734444, 0.971047664518299, 0.971047664518299, 30.22281403343689, 272.005326300932, 0.01852137464652932, 0.01852137464652931)
    10 = @("sCs", "Efna1", "Epha3", "sCs", 3, 1, 0.5101343333333334, 1.530403, 0.01907359990996641, 0.0190735999099664, 3, 1, 0.6602006666666667, 1.980602, 0.01082096921202439, 0.01082096921202439, 0.3367910269562222, 3.031119242606, 0.0002063948373882176, 0.0002063948373882175)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
